$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New match results to append starting at row 32 (A:Team1, B:Team2, C:IsTeam1Won)
$data = @(
    @("Cloud9", "Fnatic", 1),
    @("Fnatic", "Cloud9", 1),
    @("Cloud9", "Fnatic", 0),
    @("MAD Lions", "JD Gaming", 0),
    @("MAD Lions", "JD Gaming", 0),
    @("Dplus KIA", "GAM Esports", 1),
    @("GAM Esports", "Dplus KIA", 0),
    @("KT Rolster", "LNG Esports", 0),
    @("KT Rolster", "LNG Esports", 1),
    @("LNG Esports", "KT Rolster", 1),
    @("NRG", "G2 Esports", 1),
    @("G2 Esports", "NRG", 0),
    @("T1", "Bilibili Gaming", 1),
    @("Bilibili Gaming", "T1", 0),
    @("KT Rolster", "Dplus KIA", 1),
    @("Dplus KIA", "KT Rolster", 0),
    @("Fnatic", "Weibo Gaming", 1),
    @("Weibo Gaming", "Fnatic", 1),
    @("Fnatic", "Weibo Gaming", 0),
    @("G2 Esports", "Bilibili Gaming", 0),
    @("G2 Esports", "Bilibili Gaming", 1),
    @("Bilibili Gaming", "G2 Esports", 1)
)

$startRow = 32
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E52").Select()
